# Apply cryptocurrency price/volume update (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.376.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '''1.846.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("D4").Value = '''0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''240.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("D6").Value = '''0.6306'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.90%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '''0.07559'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.70%  '

$ws.Range("D9").Value = '''0.2963'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '''24.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.16%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '''2.688.08'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +45.28%  '

$ws.Range("E12").Value = '  +1.10%  '

$ws.Range("D13").Value = '''4.981'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.87%  '

$ws.Range("D14").Value = '''0.6842'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.32%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.000009974'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.50%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '''82.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.30%  '

$ws.Range("D17").Value = '''6.184'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '''29.419.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("D19").Value = '''231.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.44%  '

$ws.Range("D20").Value = '''12.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.69%  '

$ws.Range("D21").Value = '''1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").Value = '''7.582'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.45%  '

$ws.Range("D23").Value = '''1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = '''154.71'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = '''0.1391'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.66%  '

$ws.Range("D26").Value = '''8.435'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("E27").Value = '  -0.87%  '

$ws.Range("D28").Value = '''1.470'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.12%  '

$ws.Range("E29").Value = '  -3.84%  '

$ws.Range("D30").Value = '''1.266'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("D31").Value = '''4.119'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.59%  '

$ws.Range("D32").Value = '''4.016'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.47%  '

$ws.Range("D33").Value = '''1.864'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.39%  '

$ws.Range("D34").Value = '''1.158'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.90%  '

$ws.Range("B35").Value = 'RocketPoolETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D35").Value = '''2.898.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +44.32%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.7166'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '

$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").Value = '''1.249.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.13%  '

$ws.Range("E39").Value = '  -0.28%  '

$ws.Range("D40").Value = '''0.01805'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").Value = '''0.9047'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.21%  '

$ws.Range("D42").Value = '''6.076'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.70%  '

$ws.Range("D43").Value = '''0.9993'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("D44").Value = '''101.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.59%  '

$ws.Range("D45").Value = '''67.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.37%  '

$ws.Range("D46").Value = '''7.306'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("D47").Value = '''9.154'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.77%  '

$ws.Range("D48").Value = '''0.4009'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("E49").Value = '  +2.06%  '

$ws.Range("D50").Value = '''0.1121'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.60%  '

$ws.Range("D51").Value = '''0.05744'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.02%  '
